$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.43
$ws.Range("E4").Value = 16.379
$ws.Range("C6").Value = -12.998
$ws.Range("B7").Value = 5.811000000000001
$ws.Range("D7").Value = -7.858
$ws.Range("A9").Value = -21.921
$ws.Range("D10").Value = -8.415000000000001
$ws.Range("E11").Value = 16.371
$ws.Range("B12").Value = 5.568
$ws.Range("A13").Value = -22.266
$ws.Range("D13").Value = -8.1
$ws.Range("E13").Value = 16.56
$ws.Range("B14").Value = 6.077999999999999
$ws.Range("C15").Value = -12.909
$ws.Range("A16").Value = -21.767
$ws.Range("D16").Value = -8.728
$ws.Range("A18").Value = -21.962
$ws.Range("B19").Value = 7.859999999999999
$ws.Range("A20").Value = -21.296
$ws.Range("D20").Value = -8.059999999999999
$ws.Range("D24").Value = -7.153999999999999
$ws.Range("E24").Value = 16.609
$ws.Range("E25").Value = 17.165
$ws.Range("A26").Value = -21.327
$ws.Range("B26").Value = 6.141
$ws.Range("A27").Value = -21.635
$ws.Range("B27").Value = 6.486
$ws.Range("E27").Value = 16.777
$ws.Range("C28").Value = -12.965
$ws.Range("A29").Value = -21.654
$ws.Range("B29").Value = 5.997
$ws.Range("E29").Value = 17.208
$ws.Range("D32").Value = -8.328999999999999
$ws.Range("C33").Value = -11.314
$ws.Range("A35").Value = -20.228
$ws.Range("C35").Value = -12.28
$ws.Range("E35").Value = 16.277
$ws.Range("A36").Value = -20.814
$ws.Range("B37").Value = 7.515000000000001
$ws.Range("B38").Value = 5.733
$ws.Range("C38").Value = -12.043
$ws.Range("D39").Value = -7.555
$ws.Range("E40").Value = 16.598
$ws.Range("C43").Value = -12.494
$ws.Range("C44").Value = -12.708
$ws.Range("E44").Value = 16.714
$ws.Range("A45").Value = -21.64
$ws.Range("C45").Value = -12.409
$ws.Range("B47").Value = 6.622
$ws.Range("C47").Value = -12.209
$ws.Range("D47").Value = -7.323
$ws.Range("D48").Value = -7.284999999999999
$ws.Range("E49").Value = 16.562
$ws.Range("B51").Value = 5.906
$ws.Range("C51").Value = -11.937
$ws.Range("B52").Value = 6.164
$ws.Range("D52").Value = -7.459000000000001
$ws.Range("C54").Value = -13.376
$ws.Range("A55").Value = -21.724
$ws.Range("B55").Value = 6.201000000000001
$ws.Range("D56").Value = -7.972
$ws.Range("A57").Value = -21.343
$ws.Range("C57").Value = -13.103
$ws.Range("E57").Value = 16.668
$ws.Range("C62").Value = -13.591
$ws.Range("C63").Value = -11.754
$ws.Range("C67").Value = -11.229
$ws.Range("A69").Value = -21.485
$ws.Range("B69").Value = 6.479000000000001
$ws.Range("B70").Value = 6.029999999999999
$ws.Range("C70").Value = -10.995
$ws.Range("A76").Value = -21.706
$ws.Range("B76").Value = 6.703
$ws.Range("A78").Value = -20.786
$ws.Range("E80").Value = 17.021
$ws.Range("B81").Value = 5.425
$ws.Range("C81").Value = -12.753
$ws.Range("A82").Value = -21.896
$ws.Range("A83").Value = -20.889
$ws.Range("B83").Value = 6.416999999999999
$ws.Range("D84").Value = -8.267000000000001
$ws.Range("E85").Value = 16.684
$ws.Range("C88").Value = -13.081
$ws.Range("E89").Value = 17.209
$ws.Range("A93").Value = -21.688
$ws.Range("B94").Value = 6.834000000000001
$ws.Range("C96").Value = -12.998
$ws.Range("A97").Value = -21.703
$ws.Range("C99").Value = -12.729
$ws.Range("B100").Value = 6.276000000000001
$ws.Range("D100").Value = -8.335000000000001
$ws.Range("D101").Value = -7.51
$ws.Range("E101").Value = 16.677
$ws.Range("B102").Value = 6.922
